$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "1.002", "48.00", "0.00001034") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.990.50'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '1.904.83'
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '332.65'
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4633'
$ws.Range("E7").Value = '  -1.56%  '
$ws.Range("D8").Value = '0.4063'
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("D9").Value = '48.00'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D10").Value = '0.07989'
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '21.65'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").Value = '1.903.17'
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("D14").Value = '5.925'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '7.078'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '88.88'
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("D18").Value = '0.00001034'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '0.06568'
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").Value = '17.42'
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '28.980.92'
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").Value = '2.242'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("D26").Value = '2.135.88'
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").Value = '158.51'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '19.70'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").Value = '2.097'
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").Value = '5.389'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = '118.78'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '0.9782'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = '0.09371'
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("D34").Value = '1.415'
$ws.Range("E34").Value = '  +2.68%  '
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("D36").Value = '5.298'
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("D37").Value = '0.06083'
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '0.02225'
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("D39").Value = '8.408'
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").Value = '1.161'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").Value = '0.9995'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.5798'
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1821'
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '10.11'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = '1.261'
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").Value = '2.314'
$ws.Range("E46").Value = '  +11.26%  '
$ws.Range("D47").Value = '12.07'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").Value = '0.5482'
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("D49").Value = '1.902'
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("D51").Value = '47.58'
$ws.Range("E51").Value = '  +23.07%  '

# Restore default (Normal) style on column D so the explicit text
# number format does not leave a stray style index on the cells.
$ws.Range("D2:D51").Style = "Normal"

